# Update the "nomor induk" values for both users in the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> rena / Renata Agustina
$ws.Range("B2").Value = "195912312010122890"

# Row 3 -> siva / Nadya Siva
$ws.Range("B3").Value = "198311052003101001"

# Update the active selection shown in the sheet view to B2.
$ws.Range("B2").Select()
